# Realestate Update resale numbers 2025-03-02 13:15
# Appends a new row (97) of resale-number data to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

# Force the text-like columns to stay as plain text (matching the rest of
# the sheet) instead of Excel auto-converting them to a date/time/number.
# Applying NumberFormat="@" first stops the smart-parse; resetting the
# Style back to "Normal" afterwards drops the now-unneeded explicit style
# so the cell ends up as plain, unstyled text (matching the rest of the
# sheet's data rows).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-03-02"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "13:15:46"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "09"
$ws.Cells.Item($row, 4).Style = "Normal"

# Remaining columns are plain numeric resale figures (E..T).
$ws.Cells.Item($row, 5).Value = 132013
$ws.Cells.Item($row, 6).Value = 142523
$ws.Cells.Item($row, 7).Value = 171805
$ws.Cells.Item($row, 8).Value = 159089
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 147033
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193761
$ws.Cells.Item($row, 14).Value = 115247
$ws.Cells.Item($row, 15).Value = 46739
$ws.Cells.Item($row, 16).Value = 29665
$ws.Cells.Item($row, 17).Value = 70357
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 50927
$ws.Cells.Item($row, 20).Value = -1
